$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 6870
$ws.Range("I3").Value = 7136
$ws.Range("H4").Value = 1678
$ws.Range("I4").Value = 1637
$ws.Range("I5").Value = 669
$ws.Range("I6").Value = 8345
$ws.Range("H7").Value = 25991
$ws.Range("I7").Value = 24657

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("I2").Value = 72
$ws.Range("I7").Value = 290

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("I6").Value = 113
$ws.Range("I7").Value = 284

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I2").Value = 222
$ws.Range("I3").Value = 348
$ws.Range("I7").Value = 935

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("I3").Value = 72
$ws.Range("I7").Value = 218

$ws = $wb.Worksheets.Item('New City')
$ws.Range("I3").Value = 171
$ws.Range("I6").Value = 177
$ws.Range("I7").Value = 576

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I2").Value = 194
$ws.Range("H6").Value = 170
$ws.Range("I6").Value = 176
$ws.Range("I7").Value = 772
$ws.Range("I8").Value = 1475
$ws.Range("I16").Value = 71
$ws.Range("I17").Value = 35
$ws.Range("I19").Value = 694
$ws.Range("I20").Value = 606
$ws.Range("I21").Value = 109
$ws.Range("I23").Value = 244
$ws.Range("I25").Value = 129
$ws.Range("I27").Value = 214
$ws.Range("I29").Value = 1483
$ws.Range("I33").Value = 1092
$ws.Range("I36").Value = 336
$ws.Range("I40").Value = 45
$ws.Range("I41").Value = 109
$ws.Range("I42").Value = 922
$ws.Range("I50").Value = 128
$ws.Range("I52").Value = 557
$ws.Range("I53").Value = 276
$ws.Range("I54").Value = 487
$ws.Range("I55").Value = 285
$ws.Range("I57").Value = 97
$ws.Range("I60").Value = 142
$ws.Range("H63").Value = 229
$ws.Range("I63").Value = 77
$ws.Range("I64").Value = 197
$ws.Range("I65").Value = 576
$ws.Range("I67").Value = 935
$ws.Range("I69").Value = 53
$ws.Range("I71").Value = 74
$ws.Range("I72").Value = 97
$ws.Range("I73").Value = 222
$ws.Range("I78").Value = 330
$ws.Range("I79").Value = 707
$ws.Range("I83").Value = 531
$ws.Range("I84").Value = 218
$ws.Range("I85").Value = 1103
$ws.Range("I87").Value = 61
$ws.Range("I89").Value = 290
$ws.Range("I90").Value = 318
$ws.Range("I93").Value = 141
$ws.Range("I96").Value = 284
$ws.Range("I97").Value = 220
$ws.Range("I98").Value = 180
$ws.Range("H101").Value = 25991
$ws.Range("I101").Value = 24657

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("I3").Value = 194
$ws.Range("I7").Value = 531

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I3").Value = 403
$ws.Range("I4").Value = 48
$ws.Range("I6").Value = 350
$ws.Range("I7").Value = 1092

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("I3").Value = 106
$ws.Range("I7").Value = 487

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I2").Value = 438
$ws.Range("I3").Value = 509
$ws.Range("I6").Value = 409
$ws.Range("I7").Value = 1483

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I6").Value = 221
$ws.Range("I7").Value = 694

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I3").Value = 416
$ws.Range("I4").Value = 50
$ws.Range("I6").Value = 288
$ws.Range("I7").Value = 1103

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("I2").Value = 66
$ws.Range("H4").Value = 8
$ws.Range("H7").Value = 170
$ws.Range("I7").Value = 176

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("I2").Value = 34
$ws.Range("I7").Value = 109

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I2").Value = 209
$ws.Range("I6").Value = 360
$ws.Range("I7").Value = 922

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("I4").Value = 44
$ws.Range("I7").Value = 330

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("I2").Value = 87
$ws.Range("I7").Value = 285

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("I2").Value = 66
$ws.Range("I3").Value = 87
$ws.Range("I6").Value = 70
$ws.Range("I7").Value = 244

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range("I6").Value = 19
$ws.Range("I7").Value = 53

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("I6").Value = 82
$ws.Range("I7").Value = 109

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I3").Value = 229
$ws.Range("I7").Value = 707

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("I3").Value = 56
$ws.Range("I7").Value = 197

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("I6").Value = 209
$ws.Range("I7").Value = 606

$ws = $wb.Worksheets.Item('Burnside')
$ws.Range("I2").Value = 10
$ws.Range("I7").Value = 35

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("I2").Value = 96
$ws.Range("I3").Value = 112
$ws.Range("I7").Value = 336

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("I6").Value = 59
$ws.Range("I7").Value = 141

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("I2").Value = 139
$ws.Range("I7").Value = 557

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("I6").Value = 36
$ws.Range("I7").Value = 129

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("I2").Value = 34
$ws.Range("I6").Value = 116
$ws.Range("I7").Value = 180

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("I3").Value = 28
$ws.Range("I6").Value = 38
$ws.Range("I7").Value = 128

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("I4").Value = 23
$ws.Range("I7").Value = 222

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("I6").Value = 40
$ws.Range("I7").Value = 194

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("I2").Value = 36
$ws.Range("I6").Value = 143
$ws.Range("I7").Value = 220

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I6").Value = 476
$ws.Range("I7").Value = 1475

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("I4").Value = 29
$ws.Range("I7").Value = 214

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("I3").Value = 82
$ws.Range("I6").Value = 111
$ws.Range("I7").Value = 318

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("I2").Value = 37
$ws.Range("I7").Value = 97

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("I2").Value = 50
$ws.Range("I7").Value = 142

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("I6").Value = 134
$ws.Range("I7").Value = 276

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("I6").Value = 21
$ws.Range("I7").Value = 74

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("I6").Value = 48
$ws.Range("I7").Value = 97

$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Range("I3").Value = 18
$ws.Range("I6").Value = 8
$ws.Range("I7").Value = 45

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I3").Value = 236
$ws.Range("I7").Value = 772

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range("I3").Value = 12
$ws.Range("I7").Value = 61

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("I2").Value = 14
$ws.Range("I7").Value = 71
